# 9.b.1.xlsx — add a 2020 data column (N) to the "рус,англ" sheet, mirroring
# the existing 2010-2019 layout (year header in row 4, value in row 5), then
# leave the selection where the author's session ended up (N9).
#
# Note: the source diff also tweaks the Microsoft x15ac:absPath bookkeeping
# attribute in xl/workbook.xml (an absolute folder path on the author's own
# desktop, auto-stamped by Excel on save). That is not something the Excel
# object model exposes for scripting (it isn't tied to Workbook.Path/FullName
# or any other writable property) — it is incidental, machine-local save
# metadata rather than a content edit, so it is intentionally left alone here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New year header, 2020 — same style as the existing year cells (row 4).
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2020

# New data point for 2020 — same style as the existing value cells (row 5).
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 2.1

# Clear marching-ants clipboard marker left by Copy().
$excel.CutCopyMode = $false

# Match the saved selection state from the diff.
$ws.Range("N9").Select()
